# Actualización automática 2025-11-18 17:30:09
#
# Updates the figures for HIDALGO HIDALGO PEDRO GUSTAVO across the three
# report sheets, and rolls the change through the dependent total/summary
# rows in the same way the source workbook keeps them in sync.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": PORCELANATO column for this advisor ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M22").Value = -185.33

# --- Sheet "VENTA MENSUAL": noviembre column for this advisor + total row ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F22").Value = -1489.37
$wsMensual.Range("F23").Value = 7019.93

# --- Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO row + TOTAL row ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D12").Value = 6656.88
$wsCumplimiento.Range("E12").Value = 37761.12
$wsCumplimiento.Range("F12").Value = 0.1498689720383628

$wsCumplimiento.Range("D14").Value = 7019.93
$wsCumplimiento.Range("E14").Value = 48379.54101170095
$wsCumplimiento.Range("F14").Value = 0.126714747845107
